$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the published URL to https and refresh the Date ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://fhir.cqdg.ca/StructureDefinition/QCEthnicity"
$wsMeta.Range("B8").Value = "2023-04-28T18:08:06+00:00"

# --- Elements sheet: the valueCodeableConcept slice (row 7) was folded away,
#     leaving the generic Extension.value[x] row (row 6) with its
#     slicing/binding columns cleared out. ---
$wsElem = $wb.Worksheets.Item("Elements")

# Clear out the now-stale slicing / binding-strength content on row 6
# (use a leading apostrophe so the cells stay text cells, same as the
# empty-string cells already on that row, rather than turning fully blank).
$wsElem.Range("AB6").Value = "'"
$wsElem.Range("AC6").Value = "'"
$wsElem.Range("AE6").Value = "'"

# Drop the now-redundant valueCodeableConcept slice row entirely
$wsElem.Rows.Item(7).Delete()
